$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# @@ -1385,22 +1385,22 @@
$ws.Range("H15").Value = 1254.5652
$ws.Range("I15").Value = 1254.5652
$ws.Range("K15").Value = 3763.6956
$ws.Range("M15").Value = -3594.6956
# @@ -2172,22 +2172,22 @@
$ws.Range("H31").Value = 2878
$ws.Range("I31").Value = 498.5
$ws.Range("K31").Value = 1495.5
$ws.Range("M31").Value = -1265.5
# @@ -3335,22 +3335,22 @@
$ws.Range("H54").Value = 12386.667
$ws.Range("I54").Value = 4038
$ws.Range("K54").Value = 4038
$ws.Range("M54").Value = -3552
# @@ -5245,22 +5245,22 @@
$ws.Range("H92").Value = 974.8823
$ws.Range("I92").Value = 1069.5714
$ws.Range("K92").Value = 1069.5714
$ws.Range("M92").Value = 178.4286
# @@ -6301,25 +6301,22 @@
$ws.Range("H113").Value = 3956.3333
$ws.Range("I113").Value = 3956.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3956.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -702.3332999999998
$ws.Range("N113").ClearContents()
# @@ -6993,25 +6990,25 @@
$ws.Range("H127").Value = 1695.7273
$ws.Range("I127").Value = 963.6667
$ws.Range("J127").Value = 4990
$ws.Range("K127").Value = 2891.0001
$ws.Range("L127").Value = 14970
$ws.Range("M127").Value = 2068.9999
$ws.Range("N127").Value = -24890
# @@ -7195,25 +7192,25 @@
$ws.Range("H131").Value = 1001974.4
$ws.Range("J131").Value = 1534596.1
$ws.Range("L131").Value = 4603788.300000001
$ws.Range("N131").Value = -4613868.300000001
# @@ -7348,22 +7345,22 @@
$ws.Range("H134").Value = 81956
$ws.Range("J134").Value = 81956
$ws.Range("L134").Value = 81956
$ws.Range("N134").Value = -92096
# @@ -7498,22 +7495,22 @@
$ws.Range("H137").Value = 7999.778
$ws.Range("I137").Value = 7199.6
$ws.Range("K137").Value = 21598.8
$ws.Range("M137").Value = -19048.8
# @@ -7550,25 +7547,25 @@
$ws.Range("H138").Value = 2693.1333
$ws.Range("I138").Value = 2447.92
$ws.Range("J138").Value = 3919.2
$ws.Range("K138").Value = 7343.76
$ws.Range("L138").Value = 11757.6
$ws.Range("M138").Value = -2203.76
$ws.Range("N138").Value = -22037.6

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# @@ -10719,22 +10716,22 @@
$ws.Range("H61").Value = 10981567
$ws.Range("I61").Value = 14038927
$ws.Range("K61").Value = 14038927
$ws.Range("M61").Value = -14038715
# @@ -11359,25 +11356,25 @@
$ws.Range("H74").Value = 2606.3076
$ws.Range("J74").Value = 2675.5
$ws.Range("L74").Value = 2675.5
$ws.Range("N74").Value = -4423.5
# @@ -11506,25 +11503,25 @@
$ws.Range("H77").Value = 2606.3076
$ws.Range("J77").Value = 2675.5
$ws.Range("L77").Value = 13377.5
$ws.Range("N77").Value = -22113.5
# @@ -13074,22 +13071,19 @@
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# @@ -14180,25 +14174,25 @@
$ws.Range("H132").Value = 2706165
$ws.Range("I132").Value = 3587.875
$ws.Range("J132").Value = 20002658
$ws.Range("K132").Value = 10763.625
$ws.Range("L132").Value = 60007974
$ws.Range("M132").Value = -8233.625
$ws.Range("N132").Value = -60013034
# @@ -14330,22 +14324,22 @@
$ws.Range("H135").Value = 122665.25
$ws.Range("J135").Value = 122665.25
$ws.Range("L135").Value = 122665.25
$ws.Range("N135").Value = -132805.25
# @@ -14379,22 +14373,22 @@
$ws.Range("H136").Value = 10981567
$ws.Range("I136").Value = 14038927
$ws.Range("K136").Value = 42116781
$ws.Range("M136").Value = -42114231

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# @@ -18526,25 +18520,25 @@
$ws.Range("H80").Value = 1312.75
$ws.Range("I80").Value = 787.75
$ws.Range("J80").Value = 1837.75
$ws.Range("K80").Value = 787.75
$ws.Range("L80").Value = 1837.75
$ws.Range("M80").Value = 210.25
$ws.Range("N80").Value = -3833.75
# @@ -18679,25 +18673,25 @@
$ws.Range("H83").Value = 1312.75
$ws.Range("I83").Value = 787.75
$ws.Range("J83").Value = 1837.75
$ws.Range("K83").Value = 3938.75
$ws.Range("L83").Value = 9188.75
$ws.Range("M83").Value = 1053.25
$ws.Range("N83").Value = -19172.75
# @@ -19218,25 +19212,25 @@
$ws.Range("H94").Value = 1659.875
$ws.Range("I94").Value = 1758.0454
$ws.Range("J94").Value = 1443.9
$ws.Range("K94").Value = 1758.0454
$ws.Range("L94").Value = 1443.9
$ws.Range("M94").Value = -1307.0454
$ws.Range("N94").Value = -2345.9
# @@ -19858,22 +19852,22 @@
$ws.Range("H107").Value = 3885.5
$ws.Range("I107").Value = 4341.4707
$ws.Range("K107").Value = 4341.4707
$ws.Range("M107").Value = -2421.4707
# @@ -21139,22 +21133,22 @@
$ws.Range("H134").Value = 4002300.8
$ws.Range("I134").Value = 1990.1578
$ws.Range("K134").Value = 5970.4734
$ws.Range("M134").Value = -3435.4734
# @@ -21191,22 +21185,22 @@
$ws.Range("H135").Value = 59998
$ws.Range("J135").Value = 59998
$ws.Range("L135").Value = 59998
$ws.Range("N135").Value = -70138

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# @@ -26108,25 +26102,25 @@
$ws.Range("H93").Value = 52967.54
$ws.Range("I93").Value = 50582.383
$ws.Range("J93").Value = 79999.336
$ws.Range("K93").Value = 50582.383
$ws.Range("L93").Value = 79999.336
$ws.Range("M93").Value = -48710.383
$ws.Range("N93").Value = -83743.336
# @@ -26160,25 +26154,25 @@
$ws.Range("H94").Value = 1811.5385
$ws.Range("J94").Value = 2015.8889
$ws.Range("L94").Value = 2015.8889
$ws.Range("N94").Value = -2917.8889
# @@ -28138,25 +28132,25 @@
$ws.Range("H134").Value = 144969.72
$ws.Range("I134").Value = 2248.25
$ws.Range("J134").Value = 335265
$ws.Range("K134").Value = 6744.75
$ws.Range("L134").Value = 1005795
$ws.Range("M134").Value = -4209.75
$ws.Range("N134").Value = -1010865

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# @@ -34074,25 +34068,25 @@
$ws.Range("H109").Value = 5651
$ws.Range("I109").Value = 1647.125
$ws.Range("J109").Value = 21666.5
$ws.Range("K109").Value = 4941.375
$ws.Range("L109").Value = 64999.5
$ws.Range("M109").Value = -3901.375
$ws.Range("N109").Value = -67079.5
# @@ -34279,25 +34273,25 @@
$ws.Range("H113").Value = 1406.125
$ws.Range("J113").Value = 981.4
$ws.Range("L113").Value = 2944.2
$ws.Range("N113").Value = -7284.2

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# @@ -41014,22 +41008,22 @@
$ws.Range("H107").Value = 949.4483
$ws.Range("I107").Value = 959.96
$ws.Range("K107").Value = 959.96
$ws.Range("M107").Value = 960.04
# @@ -41743,25 +41737,25 @@
$ws.Range("H122").Value = 1154.7273
$ws.Range("I122").Value = 999.3333
$ws.Range("J122").Value = 1854
$ws.Range("K122").Value = 2997.9999
$ws.Range("L122").Value = 5562
$ws.Range("M122").Value = -547.9998999999998
$ws.Range("N122").Value = -10462
# @@ -42239,22 +42233,22 @@
$ws.Range("H132").Value = 4002577.8
$ws.Range("I132").Value = 2102.4736
$ws.Range("K132").Value = 6307.4208
$ws.Range("M132").Value = -3777.4208

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# @@ -43062,22 +43056,22 @@
$ws.Range("H7").Value = 4418
$ws.Range("I7").Value = 3862.75
$ws.Range("K7").Value = 3862.75
$ws.Range("M7").Value = -3750.75
# @@ -44673,22 +44667,22 @@
$ws.Range("H40").Value = 5099.8
$ws.Range("I40").Value = 5166.3335
$ws.Range("K40").Value = 5166.3335
$ws.Range("M40").Value = -5030.3335
# @@ -44967,22 +44961,22 @@
$ws.Range("H46").Value = 2980
$ws.Range("J46").Value = 2980
$ws.Range("L46").Value = 2980
$ws.Range("N46").Value = -3356
# @@ -46746,25 +46740,25 @@
$ws.Range("H82").Value = 9907.143
$ws.Range("I82").Value = 16000
$ws.Range("J82").Value = 8891.666999999999
$ws.Range("K82").Value = 16000
$ws.Range("L82").Value = 8891.666999999999
$ws.Range("M82").Value = -15639
$ws.Range("N82").Value = -9613.666999999999
# @@ -46893,25 +46887,25 @@
$ws.Range("H85").Value = 9907.143
$ws.Range("I85").Value = 16000
$ws.Range("J85").Value = 8891.666999999999
$ws.Range("K85").Value = 16000
$ws.Range("L85").Value = 8891.666999999999
$ws.Range("M85").Value = -14752
$ws.Range("N85").Value = -11387.667
# @@ -48887,22 +48881,22 @@
$ws.Range("H126").Value = 4418
$ws.Range("I126").Value = 3862.75
$ws.Range("K126").Value = 11588.25
$ws.Range("M126").Value = -9118.25
# @@ -49178,25 +49172,25 @@
$ws.Range("H132").Value = 4803.846
$ws.Range("I132").Value = 2891.6667
$ws.Range("J132").Value = 6442.857
$ws.Range("K132").Value = 8675.000100000001
$ws.Range("L132").Value = 19328.571
$ws.Range("M132").Value = -6145.000100000001
$ws.Range("N132").Value = -24388.571
# @@ -49230,22 +49224,22 @@
$ws.Range("H133").Value = 79998.5
$ws.Range("J133").Value = 79998.5
$ws.Range("L133").Value = 79998.5
$ws.Range("N133").Value = -85058.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# @@ -53073,22 +53067,22 @@
$ws.Range("H70").Value = 72729.75
$ws.Range("J70").Value = 72729.75
$ws.Range("L70").Value = 72729.75
$ws.Range("N70").Value = -73359.75
# @@ -53217,22 +53211,22 @@
$ws.Range("H73").Value = 72729.75
$ws.Range("J73").Value = 72729.75
$ws.Range("L73").Value = 72729.75
$ws.Range("N73").Value = -74913.75
# @@ -54895,25 +54889,25 @@
$ws.Range("H107").Value = 4109.7915
$ws.Range("J107").Value = 5881.6924
$ws.Range("L107").Value = 17645.0772
$ws.Range("N107").Value = -21485.0772
# @@ -56129,25 +56123,25 @@
$ws.Range("H132").Value = 387899.28
$ws.Range("I132").Value = 2851.8948
$ws.Range("J132").Value = 1433027.9
$ws.Range("K132").Value = 8555.6844
$ws.Range("L132").Value = 4299083.699999999
$ws.Range("M132").Value = -6025.6844
$ws.Range("N132").Value = -4304143.699999999
# @@ -56328,22 +56322,22 @@
$ws.Range("H136").Value = 644779.25
$ws.Range("I136").Value = 21097.867
$ws.Range("K136").Value = 63293.601
$ws.Range("M136").Value = -60743.601
